# Scheduled market-data refresh: updates Carpenter's Bench etc. price /
# profit columns (H:N) for specific Leve rows across multiple job sheets,
# reflecting the latest pulled marketboard averages.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H113").Value = 3759.4614
$ws.Range("I113").Value = 3710
$ws.Range("J113").Value = 3774.3
$ws.Range("K113").Value = 3710
$ws.Range("L113").Value = 3774.3
$ws.Range("M113").Value = -456
$ws.Range("N113").Value = -10282.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1330.8125
$ws.Range("I2").Value = 1037.2222
$ws.Range("J2").Value = 1708.2858
$ws.Range("K2").Value = 1037.2222
$ws.Range("L2").Value = 1708.2858
$ws.Range("M2").Value = -924.2221999999999
$ws.Range("N2").Value = -1934.2858

$ws.Range("H10").Value = 133.33333
$ws.Range("I10").Value = 133.33333
$ws.Range("K10").Value = 133.33333
$ws.Range("M10").Value = 36.66667000000001

$ws.Range("H88").Value = 2499.3103
$ws.Range("I88").Value = 2585.4546
$ws.Range("J88").Value = 2228.5715
$ws.Range("K88").Value = 2585.4546
$ws.Range("L88").Value = 2228.5715
$ws.Range("M88").Value = -2179.4546
$ws.Range("N88").Value = -3040.5715

$ws.Range("H91").Value = 2499.3103
$ws.Range("I91").Value = 2585.4546
$ws.Range("J91").Value = 2228.5715
$ws.Range("K91").Value = 2585.4546
$ws.Range("L91").Value = 2228.5715
$ws.Range("M91").Value = -1181.4546
$ws.Range("N91").Value = -5036.5715

$ws.Range("H116").Value = 1330.8125
$ws.Range("I116").Value = 1037.2222
$ws.Range("J116").Value = 1708.2858
$ws.Range("K116").Value = 1037.2222
$ws.Range("L116").Value = 1708.2858
$ws.Range("M116").Value = 1256.7778
$ws.Range("N116").Value = -6296.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1330.8125
$ws.Range("I3").Value = 1037.2222
$ws.Range("J3").Value = 1708.2858
$ws.Range("K3").Value = 1037.2222
$ws.Range("L3").Value = 1708.2858
$ws.Range("M3").Value = -923.2221999999999
$ws.Range("N3").Value = -1936.2858

$ws.Range("H15").Value = 20007
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 20007
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 20007
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -20461

$ws.Range("H75").Value = 4976.125
$ws.Range("I75").Value = 2829.8572
$ws.Range("K75").Value = 2829.8572
$ws.Range("M75").Value = -1893.8572

$ws.Range("H78").Value = 4976.125
$ws.Range("I78").Value = 2829.8572
$ws.Range("K78").Value = 8489.571599999999
$ws.Range("M78").Value = -3809.571599999999

$ws.Range("H86").Value = 7220.2
$ws.Range("I86").Value = 17225.334
$ws.Range("J86").Value = 2932.2856
$ws.Range("K86").Value = 17225.334
$ws.Range("L86").Value = 2932.2856
$ws.Range("M86").Value = -16102.334
$ws.Range("N86").Value = -5178.2856

$ws.Range("H89").Value = 7220.2
$ws.Range("I89").Value = 17225.334
$ws.Range("J89").Value = 2932.2856
$ws.Range("K89").Value = 86126.67
$ws.Range("L89").Value = 14661.428
$ws.Range("M89").Value = -80510.67
$ws.Range("N89").Value = -25893.428

$ws.Range("H107").Value = 1265.8572
$ws.Range("I107").Value = 1001.55
$ws.Range("J107").Value = 1926.625
$ws.Range("K107").Value = 1001.55
$ws.Range("L107").Value = 1926.625
$ws.Range("M107").Value = 918.45
$ws.Range("N107").Value = -5766.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4453.85
$ws.Range("I62").Value = 3800
$ws.Range("J62").Value = 4488.263
$ws.Range("K62").Value = 3800
$ws.Range("L62").Value = 4488.263
$ws.Range("M62").Value = -3176
$ws.Range("N62").Value = -5736.263

$ws.Range("H65").Value = 4453.85
$ws.Range("I65").Value = 3800
$ws.Range("J65").Value = 4488.263
$ws.Range("K65").Value = 19000
$ws.Range("L65").Value = 22441.315
$ws.Range("M65").Value = -15880
$ws.Range("N65").Value = -28681.315

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 105000200
$ws.Range("I9").Value = 300.5
$ws.Range("J9").Value = 315000000
$ws.Range("K9").Value = 901.5
$ws.Range("L9").Value = 945000000
$ws.Range("M9").Value = -677.5
$ws.Range("N9").Value = -945000448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 904
$ws.Range("I17").Value = 904
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 904
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -736
$ws.Range("N17").ClearContents()

$ws.Range("H70").Value = 4277.075
$ws.Range("I70").Value = 4065.4783
$ws.Range("J70").Value = 4563.353
$ws.Range("K70").Value = 4065.4783
$ws.Range("L70").Value = 4563.353
$ws.Range("M70").Value = -3795.4783
$ws.Range("N70").Value = -5103.353

$ws.Range("H73").Value = 4277.075
$ws.Range("I73").Value = 4065.4783
$ws.Range("J73").Value = 4563.353
$ws.Range("K73").Value = 4065.4783
$ws.Range("L73").Value = 4563.353
$ws.Range("M73").Value = -3129.4783
$ws.Range("N73").Value = -6435.353

$ws.Range("H102").Value = 4091.12
$ws.Range("I102").Value = 2098.2666
$ws.Range("K102").Value = 2098.2666
$ws.Range("M102").Value = -476.2665999999999

$ws.Range("H122").Value = 1653.8334
$ws.Range("I122").Value = 1708.6666
$ws.Range("J122").Value = 1599
$ws.Range("K122").Value = 5125.9998
$ws.Range("L122").Value = 4797
$ws.Range("M122").Value = -2675.9998
$ws.Range("N122").Value = -9697

$ws.Range("H126").Value = 3116.8096
$ws.Range("I126").Value = 3010.6956
$ws.Range("J126").Value = 3245.2632
$ws.Range("K126").Value = 9032.086800000001
$ws.Range("L126").Value = 9735.7896
$ws.Range("M126").Value = -6562.086800000001
$ws.Range("N126").Value = -14675.7896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2998.4736
$ws.Range("I40").Value = 2496.75
$ws.Range("J40").Value = 3858.5715
$ws.Range("K40").Value = 2496.75
$ws.Range("L40").Value = 3858.5715
$ws.Range("M40").Value = -2360.75
$ws.Range("N40").Value = -4130.5715

$ws.Range("H132").Value = 5725.423
$ws.Range("I132").Value = 1407.6428
$ws.Range("J132").Value = 10762.833
$ws.Range("K132").Value = 4222.928400000001
$ws.Range("L132").Value = 32288.499
$ws.Range("M132").Value = -1692.928400000001
$ws.Range("N132").Value = -37348.499

$ws.Range("H136").Value = 4631.625
$ws.Range("I136").Value = 2974.5715
$ws.Range("J136").Value = 6463.1055
$ws.Range("K136").Value = 8923.7145
$ws.Range("L136").Value = 19389.3165
$ws.Range("M136").Value = -6373.7145
$ws.Range("N136").Value = -24489.3165

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 7800
$ws.Range("J39").Value = 7250
$ws.Range("L39").Value = 7250
$ws.Range("N39").Value = -8076

$ws.Range("H107").Value = 546.2143
$ws.Range("I107").Value = 481.58823
$ws.Range("J107").Value = 646.0909
$ws.Range("K107").Value = 1444.76469
$ws.Range("L107").Value = 1938.2727
$ws.Range("M107").Value = 475.23531
$ws.Range("N107").Value = -5778.2727

$ws.Range("H122").Value = 1624.186
$ws.Range("I122").Value = 1649.6875
$ws.Range("K122").Value = 4949.0625
$ws.Range("M122").Value = -2499.0625
